$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.937.20"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.467.54"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.10"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.51"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D9").Value = "2.467.11"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "2.911.19"
$ws.Range("D17").Value = "62.832.10"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "2.464.15"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.14"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("E22").Value = "  +10.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.28"
$ws.Range("E25").Value = "  +21.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.66"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "658.95"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -14.25%  "
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +4.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.45"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.74"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "151.13"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "0.0₆0309"
$ws.Range("E44").Value = "  -40.06%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.23"
$ws.Range("E46").Value = "  +5.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.24"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.66"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0513"
$ws.Range("E51").Value = "  -0.67%  "
